# Insert a new weekly record into the "Coliflor" price sheet.
# This shifts the existing rows 510-609 down to 511-610 and populates
# the newly inserted row 510 with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 510 (pushes old row 510 -> 511, ..., 609 -> 610)
$ws.Rows.Item(510).Insert()

# Populate the new row with the new weekly observation
$ws.Range("A510").Value = 4
$ws.Range("B510").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C510").Value = "Los Lagos"
$ws.Range("D510").Value = 45211
$ws.Range("E510").Value = 10
$ws.Range("F510").Value = 100112008
$ws.Range("G510").Value = "Coliflor"
$ws.Range("H510").Value = "Sin especificar"
$ws.Range("I510").Value = "Primera"
$ws.Range("J510").Value = 500
$ws.Range("K510").Value = 1500
$ws.Range("L510").Value = 1500
$ws.Range("M510").Value = 1500
$ws.Range("N510").Value = "`$/unidad"
$ws.Range("O510").Value = "Región Metropolitana"
$ws.Range("P510").Value = 1500
$ws.Range("Q510").Value = 1
$ws.Range("R510").Value = "Hortaliza"
